# New Google-Forms response row (row 33) submitted by Albert Zhang.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 43908.8344866088
$ws.Range("B33").Value = "albertzhang9000@gmail.com"
$ws.Range("C33").Value = "Yes"
$ws.Range("D33").Value = "Albert Zhang"
$ws.Range("E33").Value = "Prefrosh (if you're class of 2024)"
$ws.Range("F33").Value = "Georgia"
$ws.Range("G33").Value = "Computer Science, Economics"
$ws.Range("H33").Value = "Rap, Classical, EDM"
$ws.Range("I33").Value = "Swimming & Diving, Cross Country, Track and Field"
$ws.Range("J33").Value = "Cooking/Baking, Research, Debate, Volunteering"

# Match the look of the rows above (date formatting on A, shared font elsewhere).
$ws.Range("A32").Copy()
$ws.Range("A33").PasteSpecial(-4122)

$ws.Range("B32:J32").Copy()
$ws.Range("B33:J33").PasteSpecial(-4122)
